$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 240.85715
$ws.Range("I12").Value = 287
$ws.Range("J12").Value = 179.33333
$ws.Range("K12").Value = 287
$ws.Range("L12").Value = 179.33333
$ws.Range("M12").Value = -117
$ws.Range("N12").Value = -519.3333299999999

$ws.Range("H29").Value = 678.2
$ws.Range("I29").Value = 678.2
$ws.Range("K29").Value = 2034.6
$ws.Range("M29").Value = -1753.6

$ws.Range("H31").Value = 303
$ws.Range("I31").Value = 303
$ws.Range("K31").Value = 909
$ws.Range("M31").Value = -679

$ws.Range("H51").Value = 7475
$ws.Range("I51").Value = 5000
$ws.Range("J51").Value = 7828.5713
$ws.Range("K51").Value = 5000
$ws.Range("L51").Value = 7828.5713
$ws.Range("M51").Value = -4516
$ws.Range("N51").Value = -8796.5713

$ws.Range("H53").Value = 290.9524
$ws.Range("I53").Value = 281.6154
$ws.Range("J53").Value = 306.125
$ws.Range("K53").Value = 281.6154
$ws.Range("L53").Value = 306.125
$ws.Range("M53").Value = 355.3846
$ws.Range("N53").Value = -1580.125

$ws.Range("H88").Value = 10016.667
$ws.Range("I88").Value = 10000
$ws.Range("J88").Value = 10020
$ws.Range("K88").Value = 10000
$ws.Range("L88").Value = 10020
$ws.Range("M88").Value = -9594
$ws.Range("N88").Value = -10832

$ws.Range("H91").Value = 10016.667
$ws.Range("I91").Value = 10000
$ws.Range("J91").Value = 10020
$ws.Range("K91").Value = 10000
$ws.Range("L91").Value = 10020
$ws.Range("M91").Value = -8596
$ws.Range("N91").Value = -12828

$ws.Range("H98").Value = 2443.0962
$ws.Range("I98").Value = 1775.9
$ws.Range("J98").Value = 4667.0835
$ws.Range("K98").Value = 1775.9
$ws.Range("L98").Value = 4667.0835
$ws.Range("M98").Value = -277.9000000000001
$ws.Range("N98").Value = -7663.0835

$ws.Range("H111").Value = 1520.125
$ws.Range("I111").Value = 1265.8
$ws.Range("J111").Value = 1944
$ws.Range("K111").Value = 3797.4
$ws.Range("L111").Value = 5832
$ws.Range("M111").Value = -730.3999999999996
$ws.Range("N111").Value = -11966

$ws.Range("H113").Value = 7147.5557
$ws.Range("I113").Value = 3873.8
$ws.Range("K113").Value = 3873.8
$ws.Range("M113").Value = -619.8000000000002

$ws.Range("H122").Value = 2443.0962
$ws.Range("I122").Value = 1775.9
$ws.Range("J122").Value = 4667.0835
$ws.Range("K122").Value = 5327.700000000001
$ws.Range("L122").Value = 14001.2505
$ws.Range("M122").Value = -2877.700000000001
$ws.Range("N122").Value = -18901.2505

$ws.Range("H132").Value = 31377616
$ws.Range("I132").Value = 34622684
$ws.Range("J132").Value = 8600
$ws.Range("K132").Value = 103868052
$ws.Range("L132").Value = 25800
$ws.Range("M132").Value = -103865522
$ws.Range("N132").Value = -30860

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 691.3103599999999
$ws.Range("I2").Value = 668
$ws.Range("J2").Value = 743.1111
$ws.Range("K2").Value = 668
$ws.Range("L2").Value = 743.1111
$ws.Range("M2").Value = -555
$ws.Range("N2").Value = -969.1111

$ws.Range("H47").Value = 50000
$ws.Range("J47").Value = 50000
$ws.Range("L47").Value = 50000
$ws.Range("N47").Value = -51450

$ws.Range("H49").Value = 23500
$ws.Range("J49").Value = 23500
$ws.Range("L49").Value = 23500
$ws.Range("N49").Value = -24020

$ws.Range("H104").Value = 34500
$ws.Range("J104").Value = 34500
$ws.Range("L104").Value = 34500
$ws.Range("N104").Value = -41488

$ws.Range("H116").Value = 691.3103599999999
$ws.Range("I116").Value = 668
$ws.Range("J116").Value = 743.1111
$ws.Range("K116").Value = 668
$ws.Range("L116").Value = 743.1111
$ws.Range("M116").Value = 1626
$ws.Range("N116").Value = -5331.1111

$ws.Range("H122").Value = 1907.4642
$ws.Range("I122").Value = 995.86365
$ws.Range("J122").Value = 5250
$ws.Range("K122").Value = 2987.59095
$ws.Range("L122").Value = 15750
$ws.Range("M122").Value = -537.5909499999998
$ws.Range("N122").Value = -20650

$ws.Range("H132").Value = 2093.8333
$ws.Range("J132").Value = 4806
$ws.Range("L132").Value = 14418
$ws.Range("N132").Value = -19478

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 691.3103599999999
$ws.Range("I3").Value = 668
$ws.Range("J3").Value = 743.1111
$ws.Range("K3").Value = 668
$ws.Range("L3").Value = 743.1111
$ws.Range("M3").Value = -554
$ws.Range("N3").Value = -971.1111

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H39").Value = 13106.546
$ws.Range("I39").Value = 3380.2
$ws.Range("J39").Value = 21211.834
$ws.Range("K39").Value = 3380.2
$ws.Range("L39").Value = 21211.834
$ws.Range("M39").Value = -2989.2
$ws.Range("N39").Value = -21993.834

$ws.Range("H48").Value = 46683.668
$ws.Range("J48").Value = 46683.668
$ws.Range("L48").Value = 46683.668
$ws.Range("N48").Value = -47635.668

$ws.Range("H49").Value = 13106.546
$ws.Range("I49").Value = 3380.2
$ws.Range("J49").Value = 21211.834
$ws.Range("K49").Value = 3380.2
$ws.Range("L49").Value = 21211.834
$ws.Range("M49").Value = -3198.2
$ws.Range("N49").Value = -21575.834

$ws.Range("H58").Value = 1947.0857
$ws.Range("I58").Value = 1581.5
$ws.Range("J58").Value = 5846.6665
$ws.Range("K58").Value = 1581.5
$ws.Range("L58").Value = 5846.6665
$ws.Range("M58").Value = -1378.5
$ws.Range("N58").Value = -6252.6665

$ws.Range("H99").Value = 4699.5654
$ws.Range("I99").Value = 3709.5715
$ws.Range("J99").Value = 6239.5557
$ws.Range("K99").Value = 3709.5715
$ws.Range("L99").Value = 6239.5557
$ws.Range("M99").Value = -2211.5715
$ws.Range("N99").Value = -9235.555700000001

$ws.Range("H126").Value = 4699.5654
$ws.Range("I126").Value = 3709.5715
$ws.Range("J126").Value = 6239.5557
$ws.Range("K126").Value = 11128.7145
$ws.Range("L126").Value = 18718.6671
$ws.Range("M126").Value = -8658.7145
$ws.Range("N126").Value = -23658.6671

$ws.Range("H132").Value = 3573.8096
$ws.Range("J132").Value = 4949.5
$ws.Range("L132").Value = 14848.5
$ws.Range("N132").Value = -19908.5

$ws.Range("H136").Value = 1947.0857
$ws.Range("I136").Value = 1581.5
$ws.Range("J136").Value = 5846.6665
$ws.Range("K136").Value = 4744.5
$ws.Range("L136").Value = 17539.9995
$ws.Range("M136").Value = -2194.5
$ws.Range("N136").Value = -22639.9995

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H25").Value = 4333.3335
$ws.Range("I25").Value = 3000
$ws.Range("J25").Value = 4600
$ws.Range("K25").Value = 9000
$ws.Range("L25").Value = 13800
$ws.Range("M25").Value = -8831
$ws.Range("N25").Value = -14138

$ws.Range("H30").Value = 4333.3335
$ws.Range("I30").Value = 3000
$ws.Range("J30").Value = 4600
$ws.Range("K30").Value = 9000
$ws.Range("L30").Value = 13800
$ws.Range("M30").Value = -8898
$ws.Range("N30").Value = -14004

$ws.Range("H36").Value = 0
$ws.Range("I36").Value = 0
$ws.Range("K36").Value = 0
$ws.Range("M36").ClearContents()

$ws.Range("H80").Value = 4389.5
$ws.Range("J80").Value = 4657.222
$ws.Range("L80").Value = 13971.666
$ws.Range("N80").Value = -15843.666

$ws.Range("H83").Value = 4389.5
$ws.Range("J83").Value = 4657.222
$ws.Range("L83").Value = 41914.998
$ws.Range("N83").Value = -51274.998

$ws.Range("H113").Value = 549.3276
$ws.Range("I113").Value = 560.7
$ws.Range("J113").Value = 524.05554
$ws.Range("K113").Value = 1682.1
$ws.Range("L113").Value = 1572.16662
$ws.Range("M113").Value = 487.8999999999999
$ws.Range("N113").Value = -5912.16662

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5146.1016
$ws.Range("I70").Value = 6140
$ws.Range("J70").Value = 4865.2173
$ws.Range("K70").Value = 6140
$ws.Range("L70").Value = 4865.2173
$ws.Range("M70").Value = -5870
$ws.Range("N70").Value = -5405.2173

$ws.Range("H73").Value = 5146.1016
$ws.Range("I73").Value = 6140
$ws.Range("J73").Value = 4865.2173
$ws.Range("K73").Value = 6140
$ws.Range("L73").Value = 4865.2173
$ws.Range("M73").Value = -5204
$ws.Range("N73").Value = -6737.2173

$ws.Range("H122").Value = 2468.6775
$ws.Range("I122").Value = 1993.2
$ws.Range("J122").Value = 4449.8335
$ws.Range("K122").Value = 5979.6
$ws.Range("L122").Value = 13349.5005
$ws.Range("M122").Value = -3529.6
$ws.Range("N122").Value = -18249.5005

$ws.Range("H132").Value = 2031.1136
$ws.Range("I132").Value = 1178.4
$ws.Range("J132").Value = 5347.222
$ws.Range("K132").Value = 3535.2
$ws.Range("L132").Value = 16041.666
$ws.Range("M132").Value = -1005.2
$ws.Range("N132").Value = -21101.666

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5359.231
$ws.Range("I7").Value = 2574
$ws.Range("J7").Value = 7100
$ws.Range("K7").Value = 2574
$ws.Range("L7").Value = 7100
$ws.Range("M7").Value = -2462
$ws.Range("N7").Value = -7324

$ws.Range("H122").Value = 3515.175
$ws.Range("I122").Value = 3088.121
$ws.Range("J122").Value = 5528.4287
$ws.Range("K122").Value = 9264.363000000001
$ws.Range("L122").Value = 16585.2861
$ws.Range("M122").Value = -6814.363000000001
$ws.Range("N122").Value = -21485.2861

$ws.Range("H126").Value = 5359.231
$ws.Range("I126").Value = 2574
$ws.Range("J126").Value = 7100
$ws.Range("K126").Value = 7722
$ws.Range("L126").Value = 21300
$ws.Range("M126").Value = -5252
$ws.Range("N126").Value = -26240

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 630.5
$ws.Range("I107").Value = 570
$ws.Range("K107").Value = 1710
$ws.Range("M107").Value = 210

$ws.Range("H132").Value = 6174398.5
$ws.Range("I132").Value = 656.96155
$ws.Range("K132").Value = 1970.88465
$ws.Range("M132").Value = 559.11535

$ws.Range("H136").Value = 3708.1738
$ws.Range("I136").Value = 613.7857
$ws.Range("J136").Value = 8521.666999999999
$ws.Range("K136").Value = 1841.3571
$ws.Range("L136").Value = 25565.001
$ws.Range("M136").Value = 708.6428999999998
$ws.Range("N136").Value = -30665.001
